# TC_151 test data workbook update:
#   - Update the "Default Battery Standby" (F8) and related (G8) computed
#     values on the "Add Panels" sheet to reflect the new device totals.
#   - Move the active selection from H4:L8 to K8 (last edited cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = 0.3
$ws.Range("G8").Value = 0.612

$ws.Range("K8").Select()
